$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45; everything from row 45 down shifts to row 46+.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 (copy unchanged fields from what is now row 46,
# and set the fields that differ per the target data).
$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = Get-Date -Year 2023 -Month 6 -Day 5 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = 100112037
$ws.Range("G45").Value = "Cebollín"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 100
$ws.Range("K45").Value = 6000
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = 6000
$ws.Range("N45").Value = "$/paquete 36 unidades"
$ws.Range("O45").Value = "Provincia de Diguillín"
$ws.Range("P45").Value = 167
$ws.Range("Q45").Value = 36
$ws.Range("R45").Value = "Hortaliza"
